$wb = $excel.ActiveWorkbook

# ALC row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 724.75
$ws.Range("I19").Value = 499.5
$ws.Range("J19").Value = 950
$ws.Range("K19").Value = 499.5
$ws.Range("L19").Value = 950
$ws.Range("M19").Value = -324.5
$ws.Range("N19").Value = -1300

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 317.7143
$ws.Range("I28").Value = 333
$ws.Range("K28").Value = 333
$ws.Range("M28").Value = 152

# ALC row 41
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 487.2
$ws.Range("I41").Value = 92.5
$ws.Range("J41").Value = 585.875
$ws.Range("K41").Value = 92.5
$ws.Range("L41").Value = 585.875
$ws.Range("M41").Value = 347.5
$ws.Range("N41").Value = -1465.875

# ALC row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 599.2941
$ws.Range("I92").Value = 495.57144
$ws.Range("K92").Value = 495.57144
$ws.Range("M92").Value = 752.4285600000001

# ALC row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 742.93335
$ws.Range("I107").Value = 705.5
$ws.Range("J107").Value = 785.7143
$ws.Range("K107").Value = 705.5
$ws.Range("L107").Value = 785.7143
$ws.Range("M107").Value = 1214.5
$ws.Range("N107").Value = -4625.7143

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3158.95
$ws.Range("I113").Value = 2052.5
$ws.Range("J113").Value = 5740.6665
$ws.Range("K113").Value = 2052.5
$ws.Range("L113").Value = 5740.6665
$ws.Range("M113").Value = 1201.5
$ws.Range("N113").Value = -12248.6665

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2140.3809
$ws.Range("I137").Value = 2066.125
$ws.Range("J137").Value = 2378
$ws.Range("K137").Value = 6198.375
$ws.Range("L137").Value = 7134
$ws.Range("M137").Value = -3648.375
$ws.Range("N137").Value = -12234

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1668.1578
$ws.Range("J138").Value = 4453.364
$ws.Range("L138").Value = 13360.092
$ws.Range("N138").Value = -23640.092

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3168.3171
$ws.Range("I32").Value = 2499.2812
$ws.Range("J32").Value = 5547.1113
$ws.Range("K32").Value = 2499.2812
$ws.Range("L32").Value = 5547.1113
$ws.Range("M32").Value = -2212.2812
$ws.Range("N32").Value = -6121.1113

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3250.6858
$ws.Range("I45").Value = 3320.077
$ws.Range("K45").Value = 3320.077
$ws.Range("M45").Value = -2943.077

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1288942.6
$ws.Range("I61").Value = 2573028.2
$ws.Range("K61").Value = 2573028.2
$ws.Range("M61").Value = -2572816.2

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2100
$ws.Range("I122").Value = 2133.3333
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 6399.999899999999
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -3949.999899999999
$ws.Range("N122").Value = -10900

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1288942.6
$ws.Range("I136").Value = 2573028.2
$ws.Range("K136").Value = 7719084.600000001
$ws.Range("M136").Value = -7716534.600000001

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 11252
$ws.Range("I134").Value = 12902.4
$ws.Range("K134").Value = 38707.2
$ws.Range("M134").Value = -36172.2

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1366.6666
$ws.Range("J16").Value = 1500
$ws.Range("L16").Value = 1500
$ws.Range("N16").Value = -2074

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8550.094999999999
$ws.Range("I31").Value = 9723.535
$ws.Range("J31").Value = 3504.3
$ws.Range("K31").Value = 9723.535
$ws.Range("L31").Value = 3504.3
$ws.Range("M31").Value = -9428.535
$ws.Range("N31").Value = -4094.3

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 8550.094999999999
$ws.Range("I34").Value = 9723.535
$ws.Range("J34").Value = 3504.3
$ws.Range("K34").Value = 9723.535
$ws.Range("L34").Value = 3504.3
$ws.Range("M34").Value = -9521.535
$ws.Range("N34").Value = -3908.3

# CRP row 80
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value = 21064
$ws.Range("J80").Value = 21064
$ws.Range("L80").Value = 21064
$ws.Range("N80").Value = -23310

# CRP row 83
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H83").Value = 21064
$ws.Range("J83").Value = 21064
$ws.Range("L83").Value = 63192
$ws.Range("N83").Value = -74424

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 697
$ws.Range("I107").Value = 697
$ws.Range("K107").Value = 697
$ws.Range("M107").Value = 1223

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1366.6666
$ws.Range("J113").Value = 1500
$ws.Range("L113").Value = 1500
$ws.Range("N113").Value = -5840

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1840
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550

# CUL row 17
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1500
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1500
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 4500
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -4838

# CUL row 20
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 1677.8
$ws.Range("I20").Value = 708.5
$ws.Range("J20").Value = 5555
$ws.Range("K20").Value = 2125.5
$ws.Range("L20").Value = 16665
$ws.Range("M20").Value = -1898.5
$ws.Range("N20").Value = -17119

# CUL row 34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 935
$ws.Range("J34").Value = 929.5833
$ws.Range("L34").Value = 2788.7499
$ws.Range("N34").Value = -2956.7499

# CUL row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2672.9
$ws.Range("J39").Value = 2978.625
$ws.Range("L39").Value = 8935.875
$ws.Range("N39").Value = -9523.875

# CUL row 46
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 1503.2
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1503.2
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 4509.6
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -4691.6

# CUL row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2872.5
$ws.Range("J55").Value = 2872.5
$ws.Range("L55").Value = 8617.5
$ws.Range("N55").Value = -8971.5

# CUL row 75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1550
$ws.Range("I75").Value = 1800
$ws.Range("J75").Value = 1300
$ws.Range("K75").Value = 5400
$ws.Range("L75").Value = 3900
$ws.Range("M75").Value = -4402
$ws.Range("N75").Value = -5896

# CUL row 78
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 1550
$ws.Range("I78").Value = 1800
$ws.Range("J78").Value = 1300
$ws.Range("K78").Value = 16200
$ws.Range("L78").Value = 11700
$ws.Range("M78").Value = -11208
$ws.Range("N78").Value = -21684

# CUL row 94
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 4559
$ws.Range("I94").Value = 701
$ws.Range("J94").Value = 4944.8
$ws.Range("K94").Value = 2103
$ws.Range("L94").Value = 14834.4
$ws.Range("M94").Value = -1427
$ws.Range("N94").Value = -16186.4

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 785.42
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 785.42
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2356.26
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -12436.26

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2081
$ws.Range("I122").Value = 1993.3334
$ws.Range("J122").Value = 2168.6667
$ws.Range("K122").Value = 5980.0002
$ws.Range("L122").Value = 6506.000100000001
$ws.Range("M122").Value = -3530.0002
$ws.Range("N122").Value = -11406.0001

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2370
$ws.Range("I46").Value = 2283.3333
$ws.Range("J46").Value = 2500
$ws.Range("K46").Value = 2283.3333
$ws.Range("L46").Value = 2500
$ws.Range("M46").Value = -2095.3333
$ws.Range("N46").Value = -2876

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 200.07143
$ws.Range("I55").Value = 139.4
$ws.Range("J55").Value = 233.77777
$ws.Range("K55").Value = 139.4
$ws.Range("L55").Value = 233.77777
$ws.Range("M55").Value = 33.59999999999999
$ws.Range("N55").Value = -579.77777

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 1229234.2
$ws.Range("I122").Value = 3271008
$ws.Range("J122").Value = 4169.9
$ws.Range("K122").Value = 9813024
$ws.Range("L122").Value = 12509.7
$ws.Range("M122").Value = -9810574
$ws.Range("N122").Value = -17409.7

# WVR row 86
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

# WVR row 89
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

# WVR row 101
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 12039.8
$ws.Range("J101").Value = 12039.8
$ws.Range("L101").Value = 12039.8
$ws.Range("N101").Value = -18529.8

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1243.1875
$ws.Range("I126").Value = 777.9286
$ws.Range("K126").Value = 2333.7858
$ws.Range("M126").Value = 136.2142000000003
